# Staging.PeopleReachedValues template regeneration:
# Two new "BusinessKey" columns are added to the header table on Sheet1,
# each inserted in its correct alphabetical position among the existing
# "...BusinessKey" columns:
#   - "SectorBusinessKey"    inserted right after "ResultAreaBusinessKey"
#     (and therefore right before "StatusTypeBusinessKey")
#   - "SubSectorBusinessKey" inserted right after "SubOutputBusinessKey"
#     (and therefore right before "Notes")
#
# Everything to the right of each insertion point shifts one column over,
# and the sheet's used range grows from A1:W2 to A1:Y2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert "SectorBusinessKey" immediately before the current
# "StatusTypeBusinessKey" column (column S), pushing S:W to T:X.
$ws.Columns("S").Insert()
$ws.Range("S2").Value = "SectorBusinessKey"

# Insert "SubSectorBusinessKey" immediately before the current "Notes"
# column. After the first insert, "SubOutputBusinessKey" now lives in V
# and "Notes" now lives in W, so insert before W.
$ws.Columns("W").Insert()
$ws.Range("W2").Value = "SubSectorBusinessKey"
